# Update the "July" LeetCode tracker sheet with the latest two entries
# ("Move Zeroes" and "Can Place Flowers") and add an edge-case note plus a
# hyperlink to the LeetCode submission for "Can Place Flowers".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July")

# --- Row 4: Move Zeroes -------------------------------------------------
$ws.Range("A4").Value = "easy"
$ws.Range("B4").Value = 283
$ws.Range("D4").Value = "two pointer problem"

# --- Row 5: Can Place Flowers -------------------------------------------
$ws.Range("A5").Value = "easy "
$ws.Range("B5").Value = 605

# Fill in the remaining "notes" columns after the Tag/Problem# columns so
# the new shared-string table entries are created in the same order as in
# the source workbook (easy , move zeroes , can place flowers , edge cases
# to note , <link>).
$ws.Range("C4").Value = "move zeroes "
$ws.Range("C5").Value = "can place flowers "
$ws.Range("D5").Value = "edge cases to note "

$url = "https://leetcode.com/problems/can-place-flowers/submissions/1316656978?envType=study-plan-v2&envId=leetcode-75"
$ws.Hyperlinks.Add($ws.Range("E5"), $url) | Out-Null

# Leave the selection where the author left it.
$ws.Range("B4").Select() | Out-Null
